$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 303 (pushes existing row 303..410 down to 304..411)
$ws.Rows.Item(303).Insert()

# Populate the newly inserted row 303 with the new record
$ws.Cells.Item(303, 1).Value = 3
$ws.Cells.Item(303, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(303, 3).Value = "Coquimbo"
$ws.Cells.Item(303, 4).Value = 44524
$ws.Cells.Item(303, 5).Value = 5
$ws.Cells.Item(303, 6).Value = 100112024
$ws.Cells.Item(303, 7).Value = "Choclo"
$ws.Cells.Item(303, 8).Value = "Dulce o Americano"
$ws.Cells.Item(303, 9).Value = "Primera"
$ws.Cells.Item(303, 10).Value = 90
$ws.Cells.Item(303, 11).Value = 23000
$ws.Cells.Item(303, 12).Value = 24000
$ws.Cells.Item(303, 13).Value = 23500
$ws.Cells.Item(303, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(303, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(303, 16).Value = 336
$ws.Cells.Item(303, 17).Value = 70
$ws.Cells.Item(303, 18).Value = "Hortaliza"

# Match the date-column style (s="2") used by column D elsewhere in the table
$ws.Cells.Item(303, 4).NumberFormat = $ws.Cells.Item(304, 4).NumberFormat
